$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update recalculated financial figures (rows 2-6) ---
# Row 2
$ws.Cells.Item(2, 4).Value = 36199
$ws.Cells.Item(2, 5).Value = 955
$ws.Cells.Item(2, 6).Value = 727
$ws.Cells.Item(2, 7).Value = 262
$ws.Cells.Item(2, 8).Value = 206
$ws.Cells.Item(2, 9).Value = 102
$ws.Cells.Item(2, 10).Value = 104
$ws.Cells.Item(2, 11).Value = 35784
$ws.Cells.Item(2, 12).Value = 26576
$ws.Cells.Item(2, 13).Value = 9208
$ws.Cells.Item(2, 14).Value = 6390
$ws.Cells.Item(2, 15).Value = 2818
$ws.Cells.Item(2, 16).Value = 657
$ws.Cells.Item(2, 17).Value = -257
$ws.Cells.Item(2, 18).Value = 207
$ws.Cells.Item(2, 19).Value = -823
$ws.Cells.Item(2, 20).Value = 742
$ws.Cells.Item(2, 21).Value = -999
$ws.Cells.Item(2, 22).Value = 14113
$ws.Cells.Item(2, 23).Value = 2.64
$ws.Cells.Item(2, 24).Value = 0.57
$ws.Cells.Item(2, 25).Value = 1.58
$ws.Cells.Item(2, 26).Value = 0.5600000000000001
$ws.Cells.Item(2, 27).Value = 288.63
$ws.Cells.Item(2, 28).Value = 1406.87
$ws.Cells.Item(2, 29).Value = 775
$ws.Cells.Item(2, 30).Value = 27.49
$ws.Cells.Item(2, 31).Value = 48642
$ws.Cells.Item(2, 32).Value = 0.44
$ws.Cells.Item(2, 33).Value = 500
$ws.Cells.Item(2, 34).Value = 2.35
$ws.Cells.Item(2, 35).Value = 65.05
$ws.Cells.Item(2, 36).Value = 12061185

# Row 3
$ws.Cells.Item(3, 4).Value = 35908
$ws.Cells.Item(3, 5).Value = 290
$ws.Cells.Item(3, 6).Value = 414
$ws.Cells.Item(3, 7).Value = -372
$ws.Cells.Item(3, 8).Value = -758
$ws.Cells.Item(3, 9).Value = -640
$ws.Cells.Item(3, 10).Value = -117
$ws.Cells.Item(3, 11).Value = 32918
$ws.Cells.Item(3, 12).Value = 24863
$ws.Cells.Item(3, 13).Value = 8055
$ws.Cells.Item(3, 14).Value = 5709
$ws.Cells.Item(3, 15).Value = 2346
$ws.Cells.Item(3, 16).Value = 657
$ws.Cells.Item(3, 17).Value = 2259
$ws.Cells.Item(3, 18).Value = 68
$ws.Cells.Item(3, 19).Value = -1828
$ws.Cells.Item(3, 20).Value = 950
$ws.Cells.Item(3, 21).Value = 1309
$ws.Cells.Item(3, 22).Value = 12377
$ws.Cells.Item(3, 23).Value = 0.8100000000000001
$ws.Cells.Item(3, 24).Value = -2.11
$ws.Cells.Item(3, 25).Value = -10.58
$ws.Cells.Item(3, 26).Value = -2.21
$ws.Cells.Item(3, 27).Value = 308.67
$ws.Cells.Item(3, 28).Value = 846.2
$ws.Cells.Item(3, 29).Value = -4873
$ws.Cells.Item(3, 30).Value = -15.8
$ws.Cells.Item(3, 31).Value = 43459
$ws.Cells.Item(3, 32).Value = 1.77
$ws.Cells.Item(3, 33).Value = 500
$ws.Cells.Item(3, 34).Value = 0.65
$ws.Cells.Item(3, 35).Value = -10.34
$ws.Cells.Item(3, 36).Value = 12061185

# Row 4
$ws.Cells.Item(4, 4).Value = 39369
$ws.Cells.Item(4, 5).Value = 1342
$ws.Cells.Item(4, 6).Value = 1369
$ws.Cells.Item(4, 7).Value = 383
$ws.Cells.Item(4, 8).Value = 284
$ws.Cells.Item(4, 9).Value = 205
$ws.Cells.Item(4, 10).Value = 80
$ws.Cells.Item(4, 11).Value = 33336
$ws.Cells.Item(4, 12).Value = 25023
$ws.Cells.Item(4, 13).Value = 8313
$ws.Cells.Item(4, 14).Value = 6003
$ws.Cells.Item(4, 15).Value = 2310
$ws.Cells.Item(4, 16).Value = 657
$ws.Cells.Item(4, 17).Value = 1368
$ws.Cells.Item(4, 18).Value = -969
$ws.Cells.Item(4, 19).Value = -515
$ws.Cells.Item(4, 20).Value = 1026
$ws.Cells.Item(4, 21).Value = 342
$ws.Cells.Item(4, 22).Value = 12348
$ws.Cells.Item(4, 23).Value = 3.41
$ws.Cells.Item(4, 24).Value = 0.72
$ws.Cells.Item(4, 25).Value = 3.49
$ws.Cells.Item(4, 26).Value = 0.86
$ws.Cells.Item(4, 27).Value = 301.02
$ws.Cells.Item(4, 28).Value = 847.87
$ws.Cells.Item(4, 29).Value = 1557
$ws.Cells.Item(4, 30).Value = 37.63
$ws.Cells.Item(4, 31).Value = 45691
$ws.Cells.Item(4, 32).Value = 1.28
$ws.Cells.Item(4, 33).Value = 500
$ws.Cells.Item(4, 34).Value = 0.85
$ws.Cells.Item(4, 35).Value = 32.37
$ws.Cells.Item(4, 36).Value = 12061185

# Row 5
$ws.Cells.Item(5, 4).Value = 44058
$ws.Cells.Item(5, 5).Value = 1090
$ws.Cells.Item(5, 6).Value = 1090
$ws.Cells.Item(5, 7).Value = 1138
$ws.Cells.Item(5, 8).Value = 985
$ws.Cells.Item(5, 9).Value = 862
$ws.Cells.Item(5, 10).Value = 124
$ws.Cells.Item(5, 11).Value = 36884
$ws.Cells.Item(5, 12).Value = 27854
$ws.Cells.Item(5, 13).Value = 9029
$ws.Cells.Item(5, 14).Value = 7547
$ws.Cells.Item(5, 15).Value = 1489
$ws.Cells.Item(5, 16).Value = 657
$ws.Cells.Item(5, 17).Value = 562
$ws.Cells.Item(5, 18).Value = -550
$ws.Cells.Item(5, 19).Value = 377
$ws.Cells.Item(5, 20).Value = 626
$ws.Cells.Item(5, 21).Value = -64
$ws.Cells.Item(5, 22).Value = 14129
$ws.Cells.Item(5, 23).Value = 2.47
$ws.Cells.Item(5, 24).Value = 2.24
$ws.Cells.Item(5, 25).Value = 12.73
$ws.Cells.Item(5, 26).Value = 2.81
$ws.Cells.Item(5, 27).Value = 308.48
$ws.Cells.Item(5, 28).Value = 961.5599999999999
$ws.Cells.Item(5, 29).Value = 6562
$ws.Cells.Item(5, 30).Value = 9.359999999999999
$ws.Cells.Item(5, 31).Value = 57445
$ws.Cells.Item(5, 32).Value = 1.07
$ws.Cells.Item(5, 33).Value = 500
$ws.Cells.Item(5, 34).Value = 0.8100000000000001
$ws.Cells.Item(5, 35).Value = 7.68
$ws.Cells.Item(5, 36).Value = 12061185

# Row 6
$ws.Cells.Item(6, 4).Value = 43245
$ws.Cells.Item(6, 5).Value = 1115
$ws.Cells.Item(6, 6).Value = 1115
$ws.Cells.Item(6, 7).Value = 183
$ws.Cells.Item(6, 8).Value = 16
$ws.Cells.Item(6, 9).Value = 15
$ws.Cells.Item(6, 11).Value = 34804
$ws.Cells.Item(6, 12).Value = 26450
$ws.Cells.Item(6, 13).Value = 8354
$ws.Cells.Item(6, 14).Value = 7139
$ws.Cells.Item(6, 16).Value = 685
$ws.Cells.Item(6, 17).Value = 273
$ws.Cells.Item(6, 18).Value = -812
$ws.Cells.Item(6, 19).Value = 621
$ws.Cells.Item(6, 20).Value = 576
$ws.Cells.Item(6, 21).Value = -304
$ws.Cells.Item(6, 22).Value = 15310
$ws.Cells.Item(6, 23).Value = 2.58
$ws.Cells.Item(6, 24).Value = 0.04
$ws.Cells.Item(6, 25).Value = 0.2
$ws.Cells.Item(6, 26).Value = 0.05
$ws.Cells.Item(6, 27).Value = 316.63
$ws.Cells.Item(6, 28).Value = 912.83
$ws.Cells.Item(6, 29).Value = 111
$ws.Cells.Item(6, 30).Value = 279.43
$ws.Cells.Item(6, 31).Value = 52100
$ws.Cells.Item(6, 32).Value = 0.59
$ws.Cells.Item(6, 33).Value = 500
$ws.Cells.Item(6, 34).Value = 1.62
$ws.Cells.Item(6, 35).Value = 469.61
$ws.Cells.Item(6, 36).Value = 12626426

# --- Remove stale forecast rows 7-9 data (keep columns A-C) ---
# Row 7
$ws.Range("D7:E7").ClearContents()
$ws.Range("G7:I7").ClearContents()
$ws.Range("K7:N7").ClearContents()
$ws.Range("P7:U7").ClearContents()
$ws.Range("W7:AA7").ClearContents()
$ws.Range("AC7:AI7").ClearContents()

# Row 8
$ws.Range("D8:E8").ClearContents()
$ws.Range("G8:I8").ClearContents()
$ws.Range("K8:N8").ClearContents()
$ws.Range("P8:U8").ClearContents()
$ws.Range("W8:AA8").ClearContents()
$ws.Range("AC8:AI8").ClearContents()

# Row 9
$ws.Range("D9:E9").ClearContents()
$ws.Range("G9:I9").ClearContents()
$ws.Range("K9:N9").ClearContents()
$ws.Range("P9:U9").ClearContents()
$ws.Range("W9:AA9").ClearContents()
$ws.Range("AC9:AI9").ClearContents()

